$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.326.15'
$ws.Range('D3').Value = '1.875.90'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '0.7114'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '243.08'
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.08002'
$ws.Range('E8').Value = '  +2.79%  '
$ws.Range('D9').Value = '0.3153'
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('D10').Value = '24.99'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('D11').Value = '0.08245'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').Value = '1.902.14'
$ws.Range('E12').Value = '  +1.72%  '
$ws.Range('D13').Value = '5.245'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').Value = '94.57'
$ws.Range('E14').Value = '  +3.75%  '
$ws.Range('D15').Value = '0.7123'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '6.360'
$ws.Range('E16').Value = '  +4.43%  '
$ws.Range('D17').Value = '0.000008558'
$ws.Range('E17').Value = '  +4.20%  '
$ws.Range('D18').Value = '29.345.01'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = '244.95'
$ws.Range('E19').Value = '  +2.00%  '
$ws.Range('D20').Value = '2.155.90'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '0.1556'
$ws.Range('E25').Value = '  -2.70%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '9.050'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '162.58'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').Value = '18.52'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').Value = '1.502'
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('D30').Value = '4.417'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '4.313'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').Value = '1.186'
$ws.Range('E32').Value = '  -8.02%  '
$ws.Range('D33').Value = '0.05383'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('D34').Value = '1.941'
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('D35').Value = '0.7655'
$ws.Range('E35').Value = '  +2.76%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').Value = '2.691'
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('E38').Value = '  +0.55%  '
$ws.Range('D39').Value = '1.256.07'
$ws.Range('E39').Value = '  +2.27%  '
$ws.Range('D40').Value = '2.754'
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('D41').Value = '6.501'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').Value = '0.9203'
$ws.Range('E42').Value = '  +3.74%  '
$ws.Range('D43').Value = '112.86'
$ws.Range('E43').Value = '  +1.87%  '
$ws.Range('D44').Value = '74.18'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('D45').Value = '0.00000000133'
$ws.Range('E45').Value = '  +8.61%  '
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '2.046.95'
$ws.Range('E47').Value = '  +1.35%  '
$ws.Range('D48').Value = '0.5219'
$ws.Range('E48').Value = '  +0.55%  '
$ws.Range('D49').Value = '1.802'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '9.455'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('D51').Value = '0.4358'
$ws.Range('E51').Value = '  +0.93%  '
